# Femacal de La Calera - Poroto granado
# Insert a new weekly price-report row at row 21 (pushing the existing
# rows 21-123 down to 22-124), then populate the new row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 21; Excel copies the formatting
# (e.g. the date-format style on column D) from the row above automatically.
$ws.Rows("21:21").Insert()

$ws.Range("A21").Value = 3
$ws.Range("B21").Value = 'Femacal de La Calera'
$ws.Range("C21").Value = 'Coquimbo'
$ws.Range("D21").Value = 44561
$ws.Range("E21").Value = 5
$ws.Range("F21").Value = 100112030
$ws.Range("G21").Value = 'Poroto granado'
$ws.Range("H21").Value = 'Sin especificar'
$ws.Range("I21").Value = 'Primera'
$ws.Range("J21").Value = 73
$ws.Range("K21").Value = 25000
$ws.Range("L21").Value = 26000
$ws.Range("M21").Value = 25479
$ws.Range("N21").Value = '$/saco 25 kilos'
$ws.Range("O21").Value = 'Provincia de Petorca'
$ws.Range("P21").Value = 1019
$ws.Range("Q21").Value = 25
$ws.Range("R21").Value = 'Hortaliza'
